$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 2")
$ws.Range("J4").Value = 188.4
$ws.Range("J4").Borders.Item(10).LineStyle = 1
$ws.Range("J4").Borders.Item(10).Weight = 2
